# Update "想去人数" (attendance/interest count) figures on the "展览" (Exhibition),
# "演出" (Performance) and "全部类型" (All Types) sheets to the refreshed values
# from the latest scrape, per the gh-pages data regeneration commit.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 638
$ws.Range("F4").Value = 917
$ws.Range("F5").Value = 679
$ws.Range("F6").Value = 815
$ws.Range("F8").Value = 581
$ws.Range("F9").Value = 119
$ws.Range("F10").Value = 1174
$ws.Range("F11").Value = 602
$ws.Range("F13").Value = 482
$ws.Range("F14").Value = 158
$ws.Range("F15").Value = 175
$ws.Range("F16").Value = 321
$ws.Range("F18").Value = 77
$ws.Range("F20").Value = 47
$ws.Range("F21").Value = 548
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 581

# --- 演出 (Performance) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 211
$ws.Range("F13").Value = 49

# --- 全部类型 (All Types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 638
$ws.Range("F8").Value = 917
$ws.Range("F9").Value = 679
$ws.Range("F10").Value = 815
$ws.Range("F12").Value = 581
$ws.Range("F13").Value = 119
$ws.Range("F14").Value = 1174
$ws.Range("F15").Value = 602
$ws.Range("F19").Value = 482
$ws.Range("F21").Value = 158
$ws.Range("F22").Value = 175
$ws.Range("F24").Value = 321
$ws.Range("F26").Value = 77
$ws.Range("F27").Value = 211
$ws.Range("F32").Value = 49
$ws.Range("F33").Value = 47
$ws.Range("F34").Value = 548
$ws.Range("F35").Value = 20
$ws.Range("F36").Value = 581
